$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing TEXT storage, so
# numeric-looking strings (e.g. "0.999", "311.48") are not silently
# auto-converted to the Number type by Excel. The cells style is
# captured up front and restored afterwards so the temporary "@" text
# number-format used to force the type does not leave a lasting style
# change on the cell.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

# Plain text columns (coin name, link, volume%) never look numeric to
# Excels parser, so a direct .Value assignment is safe and keeps the
# cell on its original (default) style untouched.
function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

Set-TextValue 'D2' '42.818.11'
Set-PlainValue 'E2' '  +0.26%  '
Set-TextValue 'D3' '2.539.02'
Set-PlainValue 'E3' '  -0.22%  '
Set-TextValue 'D4' '0.999'
Set-PlainValue 'E4' '  -0.04%  '
Set-TextValue 'D5' '311.48'
Set-PlainValue 'E5' '  +0.82%  '
Set-TextValue 'D6' '100.92'
Set-PlainValue 'E6' '  +3.70%  '
Set-TextValue 'D7' '0.566'
Set-PlainValue 'E7' '  -0.93%  '
Set-PlainValue 'E8' '  +0.09%  '
Set-PlainValue 'E9' '  -1.04%  '
Set-TextValue 'D10' '35.78'
Set-PlainValue 'E10' '  +0.91%  '
Set-TextValue 'D11' '0.0805'
Set-PlainValue 'E11' '  -0.01%  '
Set-TextValue 'D12' '7.33'
Set-PlainValue 'E12' '  -0.89%  '
Set-PlainValue 'E13' '  +1.48%  '
Set-TextValue 'D14' '2.930.40'
Set-PlainValue 'E14' '  -0.16%  '
Set-TextValue 'D15' '2.561.27'
Set-PlainValue 'E15' '  -0.03%  '
Set-TextValue 'D16' '15.34'
Set-PlainValue 'E16' '  -2.46%  '
Set-TextValue 'D17' '0.817'
Set-PlainValue 'E17' '  -2.07%  '
Set-TextValue 'D18' '42.819.75'
Set-PlainValue 'E18' '  +0.18%  '
Set-TextValue 'D19' '6.75'
Set-PlainValue 'E19' '  +0.25%  '
Set-TextValue 'D20' '12.35'
Set-PlainValue 'E20' '  -0.16%  '
Set-PlainValue 'E21' '  -0.14%  '
Set-TextValue 'D22' '70.15'
Set-PlainValue 'E22' '  +1.29%  '
Set-TextValue 'D23' '243.87'
Set-PlainValue 'E23' '  -1.38%  '
Set-TextValue 'D24' '2.89'
Set-PlainValue 'E24' '  -0.75%  '
Set-TextValue 'D25' '2.03'
Set-PlainValue 'E25' '  -0.49%  '
Set-PlainValue 'E26' '  +0.01%  '
Set-TextValue 'D27' '25.49'
Set-PlainValue 'E27' '  -4.19%  '
Set-PlainValue 'E28' '  -0.01%  '
Set-PlainValue 'E29' '  +0.57%  '
Set-TextValue 'D30' '38.60'
Set-PlainValue 'E30' '  -4.21%  '
Set-PlainValue 'B31' 'Filecoin'
Set-PlainValue 'C31' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D31' '5.87'
Set-PlainValue 'E31' '  +2.63%  '
Set-PlainValue 'B32' 'Monero'
Set-PlainValue 'C32' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D32' '158.71'
Set-PlainValue 'E32' '  +0.64%  '
Set-PlainValue 'E33' '  +7.23%  '
Set-TextValue 'D34' '2.68'
Set-PlainValue 'E34' '  +2.28%  '
Set-TextValue 'D35' '0.0793'
Set-PlainValue 'E35' '  -0.02%  '
Set-TextValue 'D36' '18.18'
Set-PlainValue 'E36' '  -0.71%  '
Set-PlainValue 'E37' '  -3.35%  '
Set-PlainValue 'E38' '  -4.73%  '
Set-PlainValue 'E39' '  +0.22%  '
Set-TextValue 'D40' '0.117'
Set-PlainValue 'E40' '  -0.04%  '
Set-TextValue 'D41' '4.15'
Set-PlainValue 'E41' '  +2.73%  '
Set-TextValue 'D42' '21.86'
Set-PlainValue 'E42' '  -2.36%  '
Set-PlainValue 'E43' '  +0.15%  '
Set-TextValue 'D44' '3.31'
Set-PlainValue 'E44' '  +3.91%  '
Set-TextValue 'D45' '0.0299'
Set-PlainValue 'E45' '  +0.21%  '
Set-TextValue 'D46' '1.999.70'
Set-PlainValue 'E46' '  +0.44%  '
Set-TextValue 'D47' '9.14'
Set-PlainValue 'E47' '  +1.30%  '
Set-TextValue 'D48' '2.781.61'
Set-PlainValue 'E48' '  -0.28%  '
Set-PlainValue 'E49' '  +0.63%  '
Set-TextValue 'D50' '80.29'
Set-PlainValue 'E50' '  -0.55%  '
Set-TextValue 'D51' '72.45'
Set-PlainValue 'E51' '  -0.99%  '
